$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$statuses = @("For sale", "In-transit", "For sale", "In-transit", "For sale", "In-transit", "Sold", "Sold", "Sold", "Sold")

for ($i = 0; $i -lt $statuses.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 4).Value = $statuses[$i]
}

$ws.Range("D8:D10").Select()
